$d = $word.ActiveDocument

# Locate the due-date text in the document body.
$full = $d.Content.Text
$novIdx = $full.IndexOf("November 9")

# --- "November 9" -> "November 29" -------------------------------------
# Split the run right before the "9" (using a throwaway bookmark to pin
# the boundary), change "9" to "29", then drop the throwaway bookmark.
# This leaves "November " and "29" as two separate runs, just like Word
# does when you edit in place rather than retyping the whole run.
$nineStart = $novIdx + 9
$nineEnd = $nineStart + 1
$d.Bookmarks.Add("zzSplit1", $d.Range($nineStart, $nineStart))
$d.Range($nineStart, $nineEnd).Text = "29"
$d.Bookmarks("zzSplit1").Delete()

# Word leaves its "last edit" marker (_GoBack) right after the newly
# typed text. Re-adding the bookmark by name moves it from its old
# location (after "Server") to here.
$afterTwentyNine = $nineStart + 2
$d.Bookmarks.Add("_GoBack", $d.Range($afterTwentyNine, $afterTwentyNine))

# --- ", 2017" -> ", 2018" -----------------------------------------------
# Same split trick: change the final "7" to "8" while preserving the
# run boundary right before it. Re-locate the text since the document
# shifted after the edit above.
$full = $d.Content.Text
$yearIdx = $full.IndexOf(", 2017")
$sevenStart = $yearIdx + 5
$sevenEnd = $sevenStart + 1
$d.Bookmarks.Add("zzSplit2", $d.Range($sevenStart, $sevenStart))
$d.Range($sevenStart, $sevenEnd).Text = "8"
$d.Bookmarks("zzSplit2").Delete()
